# Feature Ideas workbook update
# - insert a new "Source" column (C) between the x-marker column (B) and the
#   text column (which shifts from C to D)
# - add new rows under "Spell Casting" and "General" with author attribution
# - re-wrap long text across the sheet

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Insert the new Source column at C (shifts old C -> D) and drop the
#    formatting that Insert() copied down from column B so the new column
#    starts out blank/unstyled like the rest of the sheet.
# ---------------------------------------------------------------------------
$ws.Columns("C").Insert()
$ws.Columns("C").Clear()

# ---------------------------------------------------------------------------
# 2) Insert the two new rows inside the "Spell Casting" block (after the old
#    row 66 / new row 66) and shift everything below down by two.
# ---------------------------------------------------------------------------
$ws.Rows("67:68").Insert()

# ---------------------------------------------------------------------------
# 3) Column widths: C narrower (Source names), D keeps the old wide width.
# ---------------------------------------------------------------------------
$ws.Columns("C").ColumnWidth = 24.42578125
$ws.Columns("D").ColumnWidth = 163

# ---------------------------------------------------------------------------
# 4) Header row: new "Source" header in C1; D1 keeps "Items".
# ---------------------------------------------------------------------------
$ws.Range("C1").Value = "Source"

# ---------------------------------------------------------------------------
# 5) Source attribution tags in column C for the "Spell Casting" and
#    "General" sections.
# ---------------------------------------------------------------------------
$ewjaxRows = @(60,61,62,63,64,65,66,72,73,74,75)
foreach ($r in $ewjaxRows) {
    $ws.Cells.Item($r, 3).Value = "ewjax"
}

$ws.Cells.Item(67, 3).Value = "Dgc2002"
$ws.Cells.Item(68, 3).Value = "Isthan"
$ws.Cells.Item(76, 3).Value = "Isthan"
$ws.Cells.Item(77, 3).Value = "Cylance"

$ws.Range("C60:C68").HorizontalAlignment = -4108
$ws.Range("C72:C77").HorizontalAlignment = -4108

# ---------------------------------------------------------------------------
# 6) New feature-request rows.
# ---------------------------------------------------------------------------
$ws.Cells.Item(67, 4).Value = "If spell is overwritten (e.g. Gift of Pure Thought overwrites C2), update spell trigger list"
$ws.Cells.Item(68, 4).Value = "add feature to clear all timers but mine"
$ws.Cells.Item(76, 4).Value = 'add different coloration for "my damage"'

$cylanceText = "Another Visual and Audio alert for Randoms would be amazing, ideally you can set what will provide a trigger. Eg. `n1. All (risk is you will get these for others nearby and not in group - but you would be able to turn this off)`n2. /ran 1000 (customizable so that if someone at Angry or Ring 8 roll does roll you can be alerted)`n3. Maybe when a /ran roll is over a value. Eg. if someone does a /ran 1000 and they achieve a number >900 (these sorts of triggers often used in Raid target races etc..."
$ws.Cells.Item(77, 4).Value = $cylanceText

# ---------------------------------------------------------------------------
# 7) Wrap text across the populated stretches of the data column (skip the
#    blank gap rows so we don't leave stray empty/styled cells behind).
#    (Applied range-by-range rather than via Union, whose WrapText setter
#    only touches the first area.)
# ---------------------------------------------------------------------------
$wrapBlocks = @("D1:D1", "D3:D17", "D20:D38", "D43:D47", "D53:D56", "D59:D68", "D71:D77")
foreach ($block in $wrapBlocks) {
    $ws.Range($block).WrapText = $true
}

# Taller row for the long Cylance comment, plus a trailing spacer row.
$ws.Rows(77).RowHeight = 65.25
$ws.Rows(81).RowHeight = 15.75

# ---------------------------------------------------------------------------
# 8) Conditional formatting ("x" highlight) now spans columns B and C.
# ---------------------------------------------------------------------------
$ws.Cells.FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("B1:C1048576"))

# ---------------------------------------------------------------------------
# 9) View state: scroll the frozen pane down near the bottom and select D77,
#    matching where the author left off editing.
# ---------------------------------------------------------------------------
$ws.Activate()
$ws.Range("D77").Select()

# Window position on screen (cosmetic, best effort).
try {
    $excel.Left = 900
    $excel.Top = 1350
} catch {
}
